# Auto-generated Excel COM-interop script
# Applies the scheduled-runner price/profit refresh to the 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H38").Value = 376.07144
$ws.Range("I38").Value = 206.18182
$ws.Range("K38").Value = 618.5454599999999
$ws.Range("M38").Value = -246.5454599999999
$ws.Range("H125").Value = 594.7143
$ws.Range("I125").Value = 512.6
$ws.Range("J125").Value = 800
$ws.Range("K125").Value = 4613.400000000001
$ws.Range("L125").Value = 7200
$ws.Range("M125").Value = -2153.400000000001
$ws.Range("N125").Value = -12120
$ws.Range("H127").Value = 58824372
$ws.Range("I127").Value = 574.25
$ws.Range("J127").Value = 111112190
$ws.Range("K127").Value = 1722.75
$ws.Range("L127").Value = 333336570
$ws.Range("M127").Value = 3237.25
$ws.Range("N127").Value = -333346490
$ws.Range("H137").Value = 2177.879
$ws.Range("I137").Value = 2180.9092
$ws.Range("J137").Value = 2171.818
$ws.Range("K137").Value = 6542.7276
$ws.Range("L137").Value = 6515.454000000001
$ws.Range("M137").Value = -3992.7276
$ws.Range("N137").Value = -11615.454

$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value = 102
$ws.Range("J4").Value = 102
$ws.Range("L4").Value = 102
$ws.Range("N4").Value = -334
$ws.Range("H61").Value = 6660.4546
$ws.Range("I61").Value = 3652.3333
$ws.Range("J61").Value = 16379
$ws.Range("K61").Value = 3652.3333
$ws.Range("L61").Value = 16379
$ws.Range("M61").Value = -3440.3333
$ws.Range("N61").Value = -16803
$ws.Range("H74").Value = 5934.8213
$ws.Range("I74").Value = 2575.7727
$ws.Range("J74").Value = 18251.334
$ws.Range("K74").Value = 2575.7727
$ws.Range("L74").Value = 18251.334
$ws.Range("M74").Value = -1701.7727
$ws.Range("N74").Value = -19999.334
$ws.Range("H77").Value = 5934.8213
$ws.Range("I77").Value = 2575.7727
$ws.Range("J77").Value = 18251.334
$ws.Range("K77").Value = 12878.8635
$ws.Range("L77").Value = 91256.67
$ws.Range("M77").Value = -8510.863499999999
$ws.Range("N77").Value = -99992.67
$ws.Range("H132").Value = 1840.1818
$ws.Range("I132").Value = 1616.2222
$ws.Range("J132").Value = 2848
$ws.Range("K132").Value = 4848.6666
$ws.Range("L132").Value = 8544
$ws.Range("M132").Value = -2318.6666
$ws.Range("N132").Value = -13604
$ws.Range("H136").Value = 6660.4546
$ws.Range("I136").Value = 3652.3333
$ws.Range("J136").Value = 16379
$ws.Range("K136").Value = 10956.9999
$ws.Range("L136").Value = 49137
$ws.Range("M136").Value = -8406.999899999999
$ws.Range("N136").Value = -54237

$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 85017.914
$ws.Range("I134").Value = 1544.5555
$ws.Range("J134").Value = 335438
$ws.Range("K134").Value = 4633.666499999999
$ws.Range("L134").Value = 1006314
$ws.Range("M134").Value = -2098.666499999999
$ws.Range("N134").Value = -1011384

$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 171.71428
$ws.Range("I22").Value = 165
$ws.Range("J22").Value = 180.66667
$ws.Range("K22").Value = 165
$ws.Range("L22").Value = 180.66667
$ws.Range("M22").Value = 185
$ws.Range("N22").Value = -880.6666700000001
$ws.Range("H31").Value = 2127.0334
$ws.Range("I31").Value = 1379.72
$ws.Range("J31").Value = 5863.6
$ws.Range("K31").Value = 1379.72
$ws.Range("L31").Value = 5863.6
$ws.Range("M31").Value = -1084.72
$ws.Range("N31").Value = -6453.6
$ws.Range("H34").Value = 2127.0334
$ws.Range("I34").Value = 1379.72
$ws.Range("J34").Value = 5863.6
$ws.Range("K34").Value = 1379.72
$ws.Range("L34").Value = 5863.6
$ws.Range("M34").Value = -1177.72
$ws.Range("N34").Value = -6267.6
$ws.Range("H58").Value = 2756558.5
$ws.Range("I58").Value = 4330125
$ws.Range("K58").Value = 4330125
$ws.Range("M58").Value = -4329922
$ws.Range("H132").Value = 2682.282
$ws.Range("I132").Value = 2483.9614
$ws.Range("J132").Value = 3078.923
$ws.Range("K132").Value = 7451.8842
$ws.Range("L132").Value = 9236.769
$ws.Range("M132").Value = -4921.8842
$ws.Range("N132").Value = -14296.769
$ws.Range("H134").Value = 2316.476
$ws.Range("I134").Value = 2093.2334
$ws.Range("J134").Value = 2874.5833
$ws.Range("K134").Value = 6279.7002
$ws.Range("L134").Value = 8623.749899999999
$ws.Range("M134").Value = -3744.7002
$ws.Range("N134").Value = -13693.7499
$ws.Range("H136").Value = 2756558.5
$ws.Range("I136").Value = 4330125
$ws.Range("K136").Value = 12990375
$ws.Range("M136").Value = -12987825

$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 16667153
$ws.Range("I5").Value = 539.8889
$ws.Range("K5").Value = 1619.6667
$ws.Range("M5").Value = -1507.6667
$ws.Range("H63").Value = 3742
$ws.Range("I63").Value = 3012
$ws.Range("J63").Value = 3924.5
$ws.Range("K63").Value = 9036
$ws.Range("L63").Value = 11773.5
$ws.Range("M63").Value = -8287
$ws.Range("N63").Value = -13271.5
$ws.Range("H64").Value = 3807.889
$ws.Range("I64").Value = 2904
$ws.Range("J64").Value = 4259.8335
$ws.Range("K64").Value = 8712
$ws.Range("L64").Value = 12779.5005
$ws.Range("M64").Value = -8442
$ws.Range("N64").Value = -13319.5005
$ws.Range("H66").Value = 3742
$ws.Range("I66").Value = 3012
$ws.Range("J66").Value = 3924.5
$ws.Range("K66").Value = 27108
$ws.Range("L66").Value = 35320.5
$ws.Range("M66").Value = -23364
$ws.Range("N66").Value = -42808.5
$ws.Range("H67").Value = 3807.889
$ws.Range("I67").Value = 2904
$ws.Range("J67").Value = 4259.8335
$ws.Range("K67").Value = 8712
$ws.Range("L67").Value = 12779.5005
$ws.Range("M67").Value = -7776
$ws.Range("N67").Value = -14651.5005
$ws.Range("H122").Value = 817.8261
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 817.8261
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 7360.4349
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -12260.4349
$ws.Range("H135").Value = 16667153
$ws.Range("I135").Value = 539.8889
$ws.Range("K135").Value = 4859.0001
$ws.Range("M135").Value = -2324.0001

$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 8815.588
$ws.Range("I132").Value = 2487.9
$ws.Range("J132").Value = 17855.143
$ws.Range("K132").Value = 7463.700000000001
$ws.Range("L132").Value = 53565.429
$ws.Range("M132").Value = -4933.700000000001
$ws.Range("N132").Value = -58625.429
$ws.Range("H135").Value = 73595
$ws.Range("J135").Value = 73595
$ws.Range("L135").Value = 73595
$ws.Range("N135").Value = -83735

$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 3533.5186
$ws.Range("I132").Value = 2788.7058
$ws.Range("J132").Value = 4799.7
$ws.Range("K132").Value = 8366.117400000001
$ws.Range("L132").Value = 14399.1
$ws.Range("M132").Value = -5836.117400000001
$ws.Range("N132").Value = -19459.1
$ws.Range("H133").Value = 44339.125
$ws.Range("J133").Value = 44339.125
$ws.Range("L133").Value = 44339.125
$ws.Range("N133").Value = -49399.125
$ws.Range("H136").Value = 5482.778
$ws.Range("I136").Value = 3556.7222
$ws.Range("J136").Value = 7408.8335
$ws.Range("K136").Value = 10670.1666
$ws.Range("L136").Value = 22226.5005
$ws.Range("M136").Value = -8120.1666
$ws.Range("N136").Value = -27326.5005

$ws = $wb.Worksheets.Item(8)
$ws.Range("H29").Value = 12000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 12000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 12000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -12580
$ws.Range("H132").Value = 3493
$ws.Range("I132").Value = 2964.1667
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 8892.500100000001
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -6362.500100000001
$ws.Range("N132").Value = -25058
$ws.Range("H136").Value = 6313.0293
$ws.Range("I136").Value = 2196.5557
$ws.Range("K136").Value = 6589.6671
$ws.Range("M136").Value = -4039.6671
